$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.896.29"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.738.30"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.64"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5015"
$ws.Range("E7").Value = "  +8.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3571"
$ws.Range("E8").Value = "  +4.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.14"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07244"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.060"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.16"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.935"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "1.730.04"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.810"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001033"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06412"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.50"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.733"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").Value = "26.950.96"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.26"
$ws.Range("E24").Value = "  +4.12%  "
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.60"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.75"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").Value = "1.974.24"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.215"
$ws.Range("E29").Value = "  +5.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.54"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.041"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09527"
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.581"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.351"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02183"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05883"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.01"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2001"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.762"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6039"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.110"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.625"
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.85"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.590"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5645"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.61"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.844"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.101"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06654"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.12%  "
